# --------------------------------------------------------------------------
# This edit reorders the per-observation records on rows 5-8, 9-10, 22-23 and
# 29-31: each of those rows ends up holding the full record (species names,
# coordinates, times, free-text comments, activity/age notes, ...) that
# a DIFFERENT row held before the edit - i.e. whole rows were swapped/rotated
# among themselves, not edited field-by-field. Rows 8, 10, 17 and 36 also pick
# up an updated Taxonsorteringsordning (column B) value of 91833.
#
# Columns Y and AA (Startdatum/Slutdatum) hold the literal text '2026-02-22' on
# every row in this edit and never change, so they are intentionally left
# untouched here - re-assigning that text through .Value would make Excel's COM
# layer auto-convert it into a real date serial, which would NOT match the
# plain-text cell that is actually expected.
# --------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (record moved here from pre-edit row 7)
$ws.Range("A5").Value = 131256691
$ws.Range("B5").Value = 57884
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("I5").ClearContents()
$ws.Range("M5").Value = "äldre spår"
$ws.Range("P5").Value = "Hyttfallet, Hyttfallet, Dlr"
$ws.Range("Q5").Value = 488667
$ws.Range("R5").Value = 6665262
$ws.Range("S5").Value = 5
$ws.Range("T5").Value = "Dalarna"
$ws.Range("U5").Value = "Ludvika"
$ws.Range("V5").Value = "Dalarna"
$ws.Range("W5").Value = "Grangärde"
$ws.Range("Z5").Value = "10:55"
$ws.Range("AB5").Value = "10:55"
$ws.Range("AC5").Value = "Ringhack på gran."
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AT5").ClearContents()
$ws.Range("AW5").Value = "Karl Ericson"
$ws.Range("AX5").Value = "Karl Ericson"
$ws.Range("AY5").ClearContents()

# Row 6 (record moved here from pre-edit row 8)
$ws.Range("A6").Value = 131260583
$ws.Range("B6").Value = 57884
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("I6").ClearContents()
$ws.Range("M6").Value = "färska spår"
$ws.Range("P6").Value = "Hyttfallet, Hyttfallet, Dlr"
$ws.Range("Q6").Value = 488834
$ws.Range("R6").Value = 6665228
$ws.Range("S6").Value = 5
$ws.Range("T6").Value = "Dalarna"
$ws.Range("U6").Value = "Ludvika"
$ws.Range("V6").Value = "Dalarna"
$ws.Range("W6").Value = "Grangärde"
$ws.Range("Z6").Value = "15:30"
$ws.Range("AB6").Value = "15:30"
$ws.Range("AC6").Value = "Ringhack på tall."
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AT6").ClearContents()
$ws.Range("AW6").Value = "Karl Ericson"
$ws.Range("AX6").Value = "Karl Ericson"
$ws.Range("AY6").ClearContents()

# Row 7 (record moved here from pre-edit row 5)
$ws.Range("A7").Value = 131257424
$ws.Range("B7").Value = 79245
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("I7").ClearContents()
$ws.Range("P7").Value = "Hyttfallet, Hyttfallet, Dlr"
$ws.Range("Q7").Value = 488876
$ws.Range("R7").Value = 6665177
$ws.Range("S7").Value = 5
$ws.Range("T7").Value = "Dalarna"
$ws.Range("U7").Value = "Ludvika"
$ws.Range("V7").Value = "Dalarna"
$ws.Range("W7").Value = "Grangärde"
$ws.Range("Z7").Value = "11:33"
$ws.Range("AB7").Value = "11:33"
$ws.Range("AC7").Value = "Gran"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AT7").ClearContents()
$ws.Range("AW7").Value = "Karl Ericson"
$ws.Range("AX7").Value = "Karl Ericson"
$ws.Range("AY7").ClearContents()
$ws.Range("M7").ClearContents()

# Row 8 (record moved here from pre-edit row 6)
$ws.Range("A8").Value = 131255793
$ws.Range("B8").Value = 91833
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 5432
$ws.Range("F8").Value = "Granticka"
$ws.Range("G8").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("P8").Value = "Hyttfallet, Hyttfallet, Dlr"
$ws.Range("Q8").Value = 488817
$ws.Range("R8").Value = 6665110
$ws.Range("S8").Value = 5
$ws.Range("T8").Value = "Dalarna"
$ws.Range("U8").Value = "Ludvika"
$ws.Range("V8").Value = "Dalarna"
$ws.Range("W8").Value = "Grangärde"
$ws.Range("Z8").Value = "09:56"
$ws.Range("AB8").Value = "09:56"
$ws.Range("AC8").Value = "Flera fruktkroppar."
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AT8").ClearContents()
$ws.Range("AW8").Value = "Karl Ericson"
$ws.Range("AX8").Value = "Karl Ericson"
$ws.Range("AY8").ClearContents()
$ws.Range("M8").ClearContents()

# Row 9 (record moved here from pre-edit row 10)
$ws.Range("A9").Value = 131256423
$ws.Range("B9").Value = 57881
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 100049
$ws.Range("F9").Value = "Spillkråka"
$ws.Range("G9").Value = "Dryocopus martius"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("I9").Value = "1"
$ws.Range("K9").Value = "adult"
$ws.Range("M9").Value = "spel/sång"
$ws.Range("P9").Value = "Hyttfallet, Hyttfallet, Dlr"
$ws.Range("Q9").Value = 488671
$ws.Range("R9").Value = 6665267
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = "Dalarna"
$ws.Range("U9").Value = "Ludvika"
$ws.Range("V9").Value = "Dalarna"
$ws.Range("W9").Value = "Grangärde"
$ws.Range("Z9").Value = "10:40"
$ws.Range("AB9").Value = "10:40"
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AT9").ClearContents()
$ws.Range("AW9").Value = "Karl Ericson"
$ws.Range("AX9").Value = "Karl Ericson"
$ws.Range("AY9").ClearContents()
$ws.Range("AC9").ClearContents()

# Row 10 (record moved here from pre-edit row 9)
$ws.Range("A10").Value = 131257188
$ws.Range("B10").Value = 91833
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 5432
$ws.Range("F10").Value = "Granticka"
$ws.Range("G10").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H10").ClearContents()
$ws.Range("I10").ClearContents()
$ws.Range("P10").Value = "Hyttfallet, Hyttfallet, Dlr"
$ws.Range("Q10").Value = 488804
$ws.Range("R10").Value = 6665288
$ws.Range("S10").Value = 5
$ws.Range("T10").Value = "Dalarna"
$ws.Range("U10").Value = "Ludvika"
$ws.Range("V10").Value = "Dalarna"
$ws.Range("W10").Value = "Grangärde"
$ws.Range("Z10").Value = "11:17"
$ws.Range("AB10").Value = "11:17"
$ws.Range("AC10").Value = "Rikligt."
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AT10").ClearContents()
$ws.Range("AW10").Value = "Karl Ericson"
$ws.Range("AX10").Value = "Karl Ericson"
$ws.Range("AY10").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("M10").ClearContents()

# Row 22 (record moved here from pre-edit row 23)
$ws.Range("A22").Value = 131257316
$ws.Range("B22").Value = 79245
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = "Garnlav"
$ws.Range("G22").Value = "Alectoria sarmentosa"
$ws.Range("H22").Value = "(Ach.) Ach."
$ws.Range("I22").ClearContents()
$ws.Range("P22").Value = "Hyttfallet, Hyttfallet, Dlr"
$ws.Range("Q22").Value = 488852
$ws.Range("R22").Value = 6665209
$ws.Range("S22").Value = 5
$ws.Range("T22").Value = "Dalarna"
$ws.Range("U22").Value = "Ludvika"
$ws.Range("V22").Value = "Dalarna"
$ws.Range("W22").Value = "Grangärde"
$ws.Range("Z22").Value = "11:28"
$ws.Range("AB22").Value = "11:28"
$ws.Range("AC22").Value = "Gran"
$ws.Range("AD22").Value = $false
$ws.Range("AE22").Value = $false
$ws.Range("AG22").Value = $false
$ws.Range("AT22").ClearContents()
$ws.Range("AW22").Value = "Karl Ericson"
$ws.Range("AX22").Value = "Karl Ericson"
$ws.Range("AY22").ClearContents()

# Row 23 (record moved here from pre-edit row 22)
$ws.Range("A23").Value = 131255771
$ws.Range("B23").Value = 81230
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 1049
$ws.Range("F23").Value = "Kortskaftad ärgspik"
$ws.Range("G23").Value = "Microcalicium ahlneri"
$ws.Range("H23").Value = "Tibell"
$ws.Range("I23").ClearContents()
$ws.Range("P23").Value = "Hyttfallet, Hyttfallet, Dlr"
$ws.Range("Q23").Value = 488818
$ws.Range("R23").Value = 6665110
$ws.Range("S23").Value = 5
$ws.Range("T23").Value = "Dalarna"
$ws.Range("U23").Value = "Ludvika"
$ws.Range("V23").Value = "Dalarna"
$ws.Range("W23").Value = "Grangärde"
$ws.Range("Z23").Value = "09:54"
$ws.Range("AB23").Value = "09:54"
$ws.Range("AD23").Value = $false
$ws.Range("AE23").Value = $false
$ws.Range("AG23").Value = $false
$ws.Range("AT23").ClearContents()
$ws.Range("AW23").Value = "Karl Ericson"
$ws.Range("AX23").Value = "Karl Ericson"
$ws.Range("AY23").ClearContents()
$ws.Range("AC23").ClearContents()

# Row 29 (record moved here from pre-edit row 31)
$ws.Range("A29").Value = 131255910
$ws.Range("B29").Value = 79245
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("I29").ClearContents()
$ws.Range("P29").Value = "Hyttfallet, Hyttfallet, Dlr"
$ws.Range("Q29").Value = 488763
$ws.Range("R29").Value = 6665157
$ws.Range("S29").Value = 5
$ws.Range("T29").Value = "Dalarna"
$ws.Range("U29").Value = "Ludvika"
$ws.Range("V29").Value = "Dalarna"
$ws.Range("W29").Value = "Grangärde"
$ws.Range("Z29").Value = "10:03"
$ws.Range("AB29").Value = "10:03"
$ws.Range("AC29").Value = "Tall."
$ws.Range("AD29").Value = $false
$ws.Range("AE29").Value = $false
$ws.Range("AG29").Value = $false
$ws.Range("AT29").ClearContents()
$ws.Range("AW29").Value = "Karl Ericson"
$ws.Range("AX29").Value = "Karl Ericson"
$ws.Range("AY29").ClearContents()

# Row 30 (record moved here from pre-edit row 29)
$ws.Range("A30").Value = 131258531
$ws.Range("B30").Value = 79245
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 6425
$ws.Range("F30").Value = "Garnlav"
$ws.Range("G30").Value = "Alectoria sarmentosa"
$ws.Range("H30").Value = "(Ach.) Ach."
$ws.Range("I30").ClearContents()
$ws.Range("P30").Value = "Hyttfallet, Hyttfallet, Dlr"
$ws.Range("Q30").Value = 488725
$ws.Range("R30").Value = 6665212
$ws.Range("S30").Value = 5
$ws.Range("T30").Value = "Dalarna"
$ws.Range("U30").Value = "Ludvika"
$ws.Range("V30").Value = "Dalarna"
$ws.Range("W30").Value = "Grangärde"
$ws.Range("Z30").Value = "13:02"
$ws.Range("AB30").Value = "13:02"
$ws.Range("AC30").Value = "Gran"
$ws.Range("AD30").Value = $false
$ws.Range("AE30").Value = $false
$ws.Range("AG30").Value = $false
$ws.Range("AT30").ClearContents()
$ws.Range("AW30").Value = "Karl Ericson"
$ws.Range("AX30").Value = "Karl Ericson"
$ws.Range("AY30").ClearContents()
$ws.Range("M30").ClearContents()

# Row 31 (record moved here from pre-edit row 30)
$ws.Range("A31").Value = 131257239
$ws.Range("B31").Value = 57884
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 100109
$ws.Range("F31").Value = "Tretåig hackspett"
$ws.Range("G31").Value = "Picoides tridactylus"
$ws.Range("H31").Value = "(Linnaeus, 1758)"
$ws.Range("I31").ClearContents()
$ws.Range("M31").Value = "färska spår"
$ws.Range("P31").Value = "Hyttfallet, Hyttfallet, Dlr"
$ws.Range("Q31").Value = 488852
$ws.Range("R31").Value = 6665286
$ws.Range("S31").Value = 5
$ws.Range("T31").Value = "Dalarna"
$ws.Range("U31").Value = "Ludvika"
$ws.Range("V31").Value = "Dalarna"
$ws.Range("W31").Value = "Grangärde"
$ws.Range("Z31").Value = "11:23"
$ws.Range("AB31").Value = "11:23"
$ws.Range("AC31").Value = "Barkfläk, hagelsalva."
$ws.Range("AD31").Value = $false
$ws.Range("AE31").Value = $false
$ws.Range("AG31").Value = $false
$ws.Range("AT31").ClearContents()
$ws.Range("AW31").Value = "Karl Ericson"
$ws.Range("AX31").Value = "Karl Ericson"
$ws.Range("AY31").ClearContents()

# Rows 17 and 36 keep all their original field values in this edit; only their
# Taxonsorteringsordning (column B) is refreshed to 91833.
$ws.Range("B17").Value = 91833
$ws.Range("B36").Value = 91833
